$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Signature proportions")

# The leading "topic"/signature-name label column (column A) is being
# dropped from the "Signature proportions" sheet; every remaining column
# (topic labels in row 1, numeric signature loadings in rows 2-20) shifts
# one place to the left, which also drops the now-unused shared strings
# ("topic", "Age", "SBS5", ... "FBI/Inv") on save.
$ws.Columns.Item(1).Delete()
